$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Casos de Uso")

# Match formatting/style of the row above (row 29)
$ws.Range("B29:I29").Copy()
$ws.Range("B30:I30").PasteSpecial(-4122)

# Add new row 30 describing CU - 26 "Modificar cuenta de usuario"
# (order of entry matches shared-string insertion order in the target file)
$ws.Range("B30").Value = "CU - 26"
$ws.Range("D30").Value = "Modificar cuenta de usuario"
$ws.Range("C30").Value = "El profesor puede modificar los datos de su cuenta de usuario."
$ws.Range("E30").Value = "vacio"
$ws.Range("F30").Value = 0
$ws.Range("G30").Value = 0
$ws.Range("H30").Value = 1

$ws.Range("C30").Select()
